$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.778.59"
$ws.Range("E2").Value = "'  +0.34%  "
$ws.Range("D3").Value = "'1.629.46"
$ws.Range("E3").Value = "'  -0.03%  "
$ws.Range("E4").Value = "'  -0.85%  "
$ws.Range("D5").Value = "'214.47"
$ws.Range("E5").Value = "'  +0.20%  "
$ws.Range("E6").Value = "'  +0.07%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "'  -0.74%  "
$ws.Range("E8").Value = "'  -0.78%  "
$ws.Range("E9").Value = "'  -0.54%  "
$ws.Range("E10").Value = "'  +1.07%  "
$ws.Range("E11").Value = "'  +1.32%  "
$ws.Range("E12").Value = "'  +0.37%  "
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.645.37"
$ws.Range("E13").Value = "'  +0.74%  "
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'1.854.30"
$ws.Range("E14").Value = "'  +0.00%  "
$ws.Range("E15").Value = "'  +0.53%  "
$ws.Range("D16").Value = "'0.0₃0761"
$ws.Range("E16").Value = "'  -0.42%  "
$ws.Range("D17").Value = "'62.81"
$ws.Range("E17").Value = "'  -0.45%  "
$ws.Range("D18").Value = "'25.770.54"
$ws.Range("E18").Value = "'  +0.20%  "
$ws.Range("D19").Value = "'0.997"
$ws.Range("E19").Value = "'  -0.68%  "
$ws.Range("D21").Value = "'191.22"
$ws.Range("E21").Value = "'  -1.22%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("E22").Value = "'  +0.03%  "
$ws.Range("E23").Value = "'  +1.34%  "
$ws.Range("E24").Value = "'  -0.77%  "
$ws.Range("E25").Value = "'  +1.51%  "
$ws.Range("D26").Value = "'142.10"
$ws.Range("E26").Value = "'  +1.48%  "
$ws.Range("E27").Value = "'  +2.80%  "
$ws.Range("E28").Value = "'  +0.67%  "
$ws.Range("D29").Value = "'15.51"
$ws.Range("E29").Value = "'  +0.23%  "
$ws.Range("E30").Value = "'  +0.36%  "
$ws.Range("D31").Value = "'0.0495"
$ws.Range("E31").Value = "'  +2.39%  "
$ws.Range("E32").Value = "'  -0.05%  "
$ws.Range("E33").Value = "'  -0.28%  "
$ws.Range("E34").Value = "'  +0.27%  "
$ws.Range("E35").Value = "'  -0.49%  "
$ws.Range("D36").Value = "'0.905"
$ws.Range("E36").Value = "'  +1.27%  "
$ws.Range("D37").Value = "'1.147.21"
$ws.Range("E37").Value = "'  +4.06%  "
$ws.Range("E38").Value = "'  -0.05%  "
$ws.Range("E39").Value = "'  -2.17%  "
$ws.Range("E40").Value = "'  +0.51%  "
$ws.Range("E41").Value = "'  -0.75%  "
$ws.Range("E42").Value = "'  -0.91%  "
$ws.Range("E43").Value = "'  +0.29%  "
$ws.Range("D44").Value = "'100.66"
$ws.Range("E44").Value = "'  +0.90%  "
$ws.Range("D45").Value = "'0.804"
$ws.Range("D46").Value = "'1.764.44"
$ws.Range("E46").Value = "'  +0.10%  "
$ws.Range("B47").Value = "'BabyDogeCoin"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.0₆0108"
$ws.Range("E47").Value = "'  -0.77%  "
$ws.Range("B48").Value = "'Aave"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'55.30"
$ws.Range("E48").Value = "'  +0.57%  "
$ws.Range("B49").Value = "'RenderToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.46"
$ws.Range("E49").Value = "'  +6.70%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0511"
$ws.Range("E50").Value = "'  +1.71%  "
$ws.Range("B51").Value = "'Mantle"
$ws.Range("C51").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.415"
$ws.Range("E51").Value = "'  -0.57%  "
